{"js": "// Resume \"Skills & Abilities\" edits:\n// 1. \"Experienced with Arduino, ... Unix Shell\" -> split into 4 runs and append \", SQL\"\n// 2. \"Background in Java and Python with numerous projects completed as well as classes\"\n//    -> \"Background in Java and Python with numerous completed projects\" (split into 2 runs)\n// 3. \"Knowledgeable about hardware and software, able to adapt to new platforms easily\"\n//    -> \"Knowledge of design and software development principles, as well as good development practices\"\n// 4. Move the \"_GoBack\" bookmark from the trailing empty paragraph to right after the\n//    \"Extensive experience troubleshooting through internet and phone\" run.\n\nconst body = context.document.body;\n\n// --- 1. \"Experienced with ...\" skill line ---\nconst skillsResults = body.search(\n  \"Experienced with Arduino, C, Processing, Windows, Linux (Ubuntu, Mint), Unix Shell\",\n  { matchCase: true }\n);\nskillsResults.load(\"items\");\nawait context.sync();\n\nif (skillsResults.items.length > 0) {\n  const skillsOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n<w:r><w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:t>Experienced with</w:t></w:r>\n<w:r><w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:t xml:space=\"preserve\">: </w:t></w:r>\n<w:r><w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:t>Arduino, C, Processing, Windows, Linux (Ubuntu, Mint), Unix Shell</w:t></w:r>\n<w:r><w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:t>, SQL</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n  skillsResults.items[0].insertOoxml(skillsOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 2. \"Background in Java and Python ...\" skill line ---\nconst backgroundResults = body.search(\n  \"Background in Java and Python with numerous projects completed as well as classes\",\n  { matchCase: true }\n);\nbackgroundResults.load(\"items\");\nawait context.sync();\n\nif (backgroundResults.items.length > 0) {\n  const backgroundOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n<w:r><w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:t xml:space=\"preserve\">Background in Java and Python </w:t></w:r>\n<w:r><w:rPr><w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr><w:t>with numerous completed projects</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n  backgroundResults.items[0].insertOoxml(backgroundOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 3. \"Knowledgeable about hardware and software ...\" skill line ---\nconst knowledgeResults = body.search(\n  \"Knowledgeable about hardware and software, able to adapt to new platforms easily\",\n  { matchCase: true }\n);\nknowledgeResults.load(\"items\");\nawait context.sync();\n\nif (knowledgeResults.items.length > 0) {\n  knowledgeResults.items[0].insertText(\n    \"Knowledge of design and software development principles, as well as good development practices\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// --- 4. Move the \"_GoBack\" bookmark next to the troubleshooting bullet ---\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst troubleshootingResults = body.search(\n  \"Extensive experience troubleshooting through internet and phone\",\n  { matchCase: true }\n);\ntroubleshootingResults.load(\"items\");\nawait context.sync();\n\nif (troubleshootingResults.items.length > 0) {\n  const endOfTroubleshooting = troubleshootingResults.items[0].getRange(Word.RangeLocation.end);\n  endOfTroubleshooting.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Resume \"Skills & Abilities\" edits:\n# 1. \"Experienced with Arduino, ... Unix Shell\" -> split into runs and append \", SQL\"\n# 2. \"Background in Java and Python with numerous projects completed as well as classes\"\n#    -> \"Background in Java and Python with numerous completed projects\"\n# 3. \"Knowledgeable about hardware and software, able to adapt to new platforms easily\"\n#    -> \"Knowledge of design and software development principles, as well as good development practices\"\n# 4. Move the \"_GoBack\" bookmark from the trailing empty paragraph to right after the\n#    \"Extensive experience troubleshooting through internet and phone\" run.\n\n$d = $word.ActiveDocument\n\n# --- 1. \"Experienced with ...\" skill line ---\n$rng1 = $d.Content\n$find1 = $rng1.Find\n$find1.ClearFormatting()\n$find1.Text = \"Experienced with Arduino, C, Processing, Windows, Linux (Ubuntu, Mint), Unix Shell\"\n$found1 = $find1.Execute()\nif ($found1) {\n    $rng1.Text = \"Experienced with\"\n    $rng1.Collapse(0)\n    $rng1.InsertAfter(\": \")\n    $rng1.Collapse(0)\n    $rng1.InsertAfter(\"Arduino, C, Processing, Windows, Linux (Ubuntu, Mint), Unix Shell\")\n    $rng1.Collapse(0)\n    $rng1.InsertAfter(\", SQL\")\n}\n\n# --- 2. \"Background in Java and Python ...\" skill line ---\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.ClearFormatting()\n$find2.Text = \"Background in Java and Python with numerous projects completed as well as classes\"\n$found2 = $find2.Execute()\nif ($found2) {\n    $rng2.Text = \"Background in Java and Python \"\n    $rng2.Collapse(0)\n    $rng2.InsertAfter(\"with numerous completed projects\")\n}\n\n# --- 3. \"Knowledgeable about hardware and software ...\" skill line ---\n$rng3 = $d.Content\n$find3 = $rng3.Find\n$find3.ClearFormatting()\n$find3.Text = \"Knowledgeable about hardware and software, able to adapt to new platforms easily\"\n$found3 = $find3.Execute()\nif ($found3) {\n    $rng3.Text = \"Knowledge of design and software development principles, as well as good development practices\"\n}\n\n# --- 4. Move the \"_GoBack\" bookmark next to the troubleshooting bullet ---\n$oldBookmark = $d.Bookmarks(\"_GoBack\")\n$oldBookmark.Delete()\n\n$rng4 = $d.Content\n$find4 = $rng4.Find\n$find4.ClearFormatting()\n$find4.Text = \"Extensive experience troubleshooting through internet and phone\"\n$found4 = $find4.Execute()\nif ($found4) {\n    $endPos = $rng4.End\n    # Temporarily insert a one-character placeholder right after the bullet text so we\n    # have a non-collapsed range to anchor the bookmark on, then delete the placeholder,\n    # leaving the bookmark collapsed right after the text (mirrors how Word relocates\n    # its auto-maintained \"_GoBack\" bookmark to the most recently edited spot).\n    $rng4.Collapse(0)\n    $rng4.InsertAfter(\"X\")\n    $placeholderRange = $d.Range($endPos, $endPos + 1)\n    $d.Bookmarks.Add(\"_GoBack\", $placeholderRange)\n    $placeholderRange2 = $d.Range($endPos, $endPos + 1)\n    $placeholderRange2.Delete()\n}\n"}
